$wb = $excel.ActiveWorkbook

# --- Sheets ---
$nrc = $wb.Worksheets.Item("NRC")
$subtask = $wb.Worksheets.Item("Subtask")

# --- Update existing data (NRC sheet) ---
# J4: 0.01 -> 4
$nrc.Range("J4").Value = 4
# J6: 108 -> 200
$nrc.Range("J6").Value = 200

# --- Add new row 9 to NRC sheet (copy row 7's formatting/types, then edit) ---
$nrc.Range("A7:K7").Copy()
$nrc.Range("A9").PasteSpecial()
$nrc.Range("A9").Value = "QZL0071"
$nrc.Range("J9").Value = 12
$nrc.Range("K9").Value = 0

# --- Selections / active sheet ---
# NRC sheet becomes the active/selected sheet with a new selection
$nrc.Activate()
$nrc.Range("L10").Select()

# Subtask sheet selection remains E1 (unchanged), but it's no longer the active tab
$subtask.Range("E1").Select()

# Re-activate NRC sheet so it is the active tab on save
$nrc.Activate()
